# Daily attendance processing - 2026-01-18 09:33:56
# Reorders the "Recorded By" (column G) values on each row: the two (or
# first two, for the triple-author rows) names/emails in the comma
# separated list are swapped into the new canonical order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "system, System, backup@backdoor.com" = "System, system, backup@backdoor.com"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value()
    if ($null -ne $val -and $val -ne "" -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
